# Update zero prandtl number
# Adds three new blocks of runs (rows 61-64, 66-69, 71-74) to the "RBC flux"
# sheet, mirroring the existing Pr=1 128x128 grid-resolution blocks found
# earlier in the sheet (e.g. rows 12-14), but under three new section
# headers describing truncation-number variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderRow($row, $text) {
    $ws.Cells.Item($row, 1).Value = $text
}

function Set-DataRow($row, $slurmNum, $raTq, $note) {
    $ws.Cells.Item($row, 1).Value = $slurmNum
    $ws.Cells.Item($row, 1).Font.Name = "Unset"

    $ws.Cells.Item($row, 2).Value = $raTq
    $ws.Cells.Item($row, 3).Value = 1
    $ws.Cells.Item($row, 4).Value = "2pi/10"
    $ws.Cells.Item($row, 5).Value = 1
    $ws.Cells.Item($row, 6).Value = 128
    $ws.Cells.Item($row, 7).Value = 128
    $ws.Cells.Item($row, 8).Value = 10
    $ws.Cells.Item($row, 9).Value = 10
    $ws.Cells.Item($row, 10).Value = 0

    $ws.Cells.Item($row, 11).Value = 0.001
    $ws.Cells.Item($row, 11).NumberFormat = "0.00E+00"

    $ws.Cells.Item($row, 12).Value = 10
    $ws.Cells.Item($row, 13).Value = 0.01

    $ws.Cells.Item($row, 14).Value = $note
    $ws.Cells.Item($row, 14).Font.Bold = $true
}

# --- Section: nx_trunc_num=1, nz_trunc_num=2 ---
Set-HeaderRow 61 "nx_trunc_num=1, nz_trunc_num=2"
Set-DataRow 62 14965608 30000 "steady tilted roll"
Set-DataRow 63 14965609 40000 "direction reversing tilted roll"
Set-DataRow 64 14965611 60000 "modulated traveling tilted roll"

# --- Section: nx_trunc_num=1, nz_trunc_num=2, with A_noise=0.01 ---
Set-HeaderRow 66 "nx_trunc_num=1, nz_trunc_num=2, with A_noise=0.01"
Set-DataRow 67 14965912 30000 "steady tilted roll"
Set-DataRow 68 14965914 40000 "direction reversing tilted roll"
Set-DataRow 69 14965915 60000 "modulated traveling tilted roll"

# --- Section: nx_trunc_num=1, nz_trunc_num=1000 ---
Set-HeaderRow 71 "nx_trunc_num=1, nz_trunc_num=1000"
Set-DataRow 72 14965918 30000 "steady tilted roll"
Set-DataRow 73 14965919 40000 "direction reversing tilted roll"
Set-DataRow 74 14965920 60000 "modulated traveling tilted roll"

# Move the view / selection the way the author last left it: scrolled down
# so row 58 is at the top, with A75 (just past the new data) selected.
$ws.Range("A75").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
